$d = $word.ActiveDocument

# 1. Merge the three runs around "proposedNRT" (which are split apart by
#    w:proofErr spell-check markers) back into a single run with the
#    combined text, by doing a Find & Replace over the whole phrase.
$oldText = ": Describe the experience of the PI and Co-PIs with leading or participating in STEM education and training over the past five years. Describe any overlap and/or complementarity between the training and the proposedNRT program."
$newText = ": Describe the experience of the PI and Co-PIs with leading or participating in STEM education and training over the past five years. Describe any overlap and/or complementarity between the training and the proposedNRT program."
$d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null

# 2. Insert a new, empty paragraph between the descriptive paragraph and
#    the final (bookmarked) paragraph.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore()

# 3. Add a "1 page" run at the start of the final paragraph (the one that
#    still holds the _GoBack bookmark).
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertBefore("1 page")
